$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the "(c) 2020 ... Creative
# Commons Attribution" paragraph that follow the "Requisitos" section at the
# end of the document, then remove that whole block (including the blank
# paragraph that precedes "Ver no Jupiter ..."), leaving only the single
# blank paragraph that originally trailed the copyright notice.

$startPara = $null
$endPara = $null
$count = $d.Paragraphs.Count

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($startPara -eq $null -and $t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $p
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $deleteStart = $startPara.Previous().Range.Start
    $deleteEnd = $endPara.Range.End
    $range = $d.Range($deleteStart, $deleteEnd)
    $range.Delete()
}
